$d = $word.ActiveDocument

# The document originally has 5 paragraphs:
#   1. "This is my code repository!"
#   2. "This is version 1"              (has a pPr/rPr/rFonts eastAsia hint)
#   3. ""                                (empty, pPr/rPr/rFonts eastAsia hint)
#   4. ""                                (empty, pPr/rPr/rFonts eastAsia hint)
#   5. "This is version 2" + bookmarkStart/bookmarkEnd ("_GoBack")
#
# The revert collapses this back down to 2 paragraphs:
#   1. "This is my code repository!"
#   2. "This is version 1" + bookmarkStart/bookmarkEnd ("_GoBack")
# i.e. the two blank paragraphs and the "This is version 2" text are removed,
# while the _GoBack bookmark around the old paragraph 5 is preserved and ends
# up trailing paragraph 2.

# Remove the two blank paragraphs (paragraph 3, twice, since paragraph 4
# shifts into index 3 once paragraph 3 is gone).
$d.Paragraphs(3).Range.Delete()
$d.Paragraphs(3).Range.Delete()

# Clear the text of the former "This is version 2" paragraph (now paragraph 3)
# but leave its bookmarkStart/bookmarkEnd marks in place.
$p3 = $d.Paragraphs(3)
$r = $p3.Range
$r.End = $r.End - 0
$r.Text = ""

# Merge the now-empty paragraph 3 (holding only the bookmark) up into
# paragraph 2 by deleting the paragraph mark that separates them.
$p2 = $d.Paragraphs(2)
$markRange = $d.Range($p2.Range.End - 1, $p2.Range.End)
$markRange.Delete()
